$d = $word.ActiveDocument
$t = $d.Tables.Item(1)
$nl = [char]11

# New lattice-multiplication exercise data, in reading order (row-major),
# matching the existing 5 x 3 grid of cells one-for-one.
$cells = @(
    @("29 x 69", "  6    9", "2|    |", "9|    |"),
    @("15 x 51", "  5    1", "1|    |", "5|    |"),
    @("73 x 28", "  2    8", "7|    |", "3|    |"),
    @("34 x 81", "  8    1", "3|    |", "4|    |"),
    @("14 x 96", "  9    6", "1|    |", "4|    |"),
    @("58 x 40", "  4    0", "5|    |", "8|    |"),
    @("10 x 66", "  6    6", "1|    |", "0|    |"),
    @("82 x 86", "  8    6", "8|    |", "2|    |"),
    @("28 x 53", "  5    3", "2|    |", "8|    |"),
    @("12 x 72", "  7    2", "1|    |", "2|    |"),
    @("53 x 97", "  9    7", "5|    |", "3|    |"),
    @("43 x 11", "  1    1", "4|    |", "3|    |"),
    @("96 x 82", "  8    2", "9|    |", "6|    |"),
    @("51 x 44", "  4    4", "5|    |", "1|    |"),
    @("78 x 98", "  9    8", "7|    |", "8|    |")
)

$rows = 5
$cols = 3
for ($i = 0; $i -lt $cells.Count; $i++) {
    $r = [math]::Floor($i / $cols) + 1
    $c = ($i % $cols) + 1
    $vals = $cells[$i]
    $newText = $vals[0] + $nl + $vals[1] + $nl + "  ----" + $nl + $vals[2] + $nl + $vals[3]
    $t.Cell($r, $c).Range.Text = $newText
}
